# Generate Report for Handoff
# Updates the localization status Overview / zh-cn / de-de sheets:
#  - Status moves from "In Translation" to "Ready for handoff"
#  - The "generate date" / "handoff datetime" timestamps are refreshed
#  - The now-wider "Ready for handoff" text needs wider status columns

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns (E2, F2) + generate date (G2)
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-14 00:52:48"

# zh-cn sheet: Status column (C2) + Latest Handoff Datetime (H2)
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-14 00:52:40"

# de-de sheet: Status column (C2) + Latest Handoff Datetime (H2)
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-14 00:52:48"

# Widen the Status columns so the longer "Ready for handoff" text fits.
# (The engine snaps ColumnWidth to an integer-pixel grid on save, so 16.3
# is the input that lands on the nearest achievable stored width to the
# target 17.2159881591797, i.e. 17.166666666666668.)
$overview.Columns.Item(5).ColumnWidth = 16.3
$overview.Columns.Item(6).ColumnWidth = 16.3
$zhcn.Columns.Item(3).ColumnWidth = 16.3
$dede.Columns.Item(3).ColumnWidth = 16.3
